$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the auto-updating "datetimeFigureOut" date field cached on the
#    slide master and every slide layout (6/18/2021 -> 6/24/2021).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "6/18/2021") {
                $shp.TextFrame.TextRange.Text = "6/24/2021"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 2: split the "FFR, SOFR Qual Forecasts" label into
#    "FFR, " + "SOFR Baseline " + "Forecasts" (new baseline qual forecasts).
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(2)
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.TextRange.Text -eq "FFR, SOFR Qual Forecasts") {
        $target = $candidate
        break
    }
}

if ($target -ne $null) {
    $tr = $target.TextFrame.TextRange
    # "FFR, SOFR Qual Forecasts"
    #  1234567890123456789012345  (1-based)
    # chars 6..15 = "SOFR Qual " -> "SOFR Baseline "
    $mid = $tr.Characters(6, 10)
    $mid.Text = "SOFR Baseline "
}
